$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 27.43619555144485
$ws.Range("C2").Value = 23.42188548628871
$ws.Range("D2").Value = 0.2638380750548214
$ws.Range("E2").Value = 0.267513168418351
$ws.Range("F2").Value = 1.303473136068624
$ws.Range("G2").Value = 0.2636843208182839
$ws.Range("H2").Value = 0.2641790516651121
$ws.Range("B3").Value = 26.39278025933316
$ws.Range("C3").Value = 22.54125019845502
$ws.Range("D3").Value = 0.2641141330047871
$ws.Range("E3").Value = 0.2676019174960007
$ws.Range("F3").Value = 1.302682128242274
$ws.Range("G3").Value = 0.2639668468031766
$ws.Range("H3").Value = 0.264440143872206
$ws.Range("B4").Value = 25.28917663327386
$ws.Range("C4").Value = 21.63136539485483
$ws.Range("D4").Value = 0.2644153266869225
$ws.Range("E4").Value = 0.2676933314007488
$ws.Range("F4").Value = 0.2645476000079352
$ws.Range("G4").Value = 0.2642747230194855
$ws.Range("H4").Value = 0.264725874812176
$ws.Range("B5").Value = 25.28376857992589
$ws.Range("C5").Value = 21.626797444188
$ws.Range("D5").Value = 0.264416751738639
$ws.Range("E5").Value = 0.2676937896096514
$ws.Range("F5").Value = 0.2645490037671267
$ws.Range("G5").Value = 0.2642761815349058
$ws.Range("H5").Value = 0.2647272224369829
$ws.Range("B6").Value = 24.59202864437152
$ws.Range("C6").Value = 21.06759092294143
$ws.Range("D6").Value = 0.2645988760186596
$ws.Range("E6").Value = 1.292630545529656
$ws.Range("F6").Value = 0.2647268120882884
$ws.Range("G6").Value = 0.2644625844995955
$ws.Range("H6").Value = 0.2648994469340974
$ws.Range("B7").Value = 25.15792493393231
$ws.Range("C7").Value = 21.49851741625266
$ws.Range("D7").Value = 0.2644400093939182
$ws.Range("E7").Value = 0.2677066543105999
$ws.Range("F7").Value = 0.2645713249478115
$ws.Range("G7").Value = 0.2643003681189386
$ws.Range("H7").Value = 1.301227227379357
$ws.Range("B8").Value = 26.45907828948702
$ws.Range("C8").Value = 22.59721805105877
$ws.Range("D8").Value = 0.2640966130629014
$ws.Range("E8").Value = 0.2675962851563812
$ws.Range("F8").Value = 1.302732246244567
$ws.Range("G8").Value = 0.2639489161091094
$ws.Range("H8").Value = 0.2644235878006646
$ws.Range("B9").Value = 30.04102059701833
$ws.Range("C9").Value = 25.61848406093264
$ws.Range("D9").Value = 0.2631458688307897
$ws.Range("E9").Value = 0.26729062467343
$ws.Range("F9").Value = 1.305466218960807
$ws.Range("G9").Value = 0.2629759343869661
$ws.Range("H9").Value = 0.2635242840646136
$ws.Range("B10").Value = 42.47587748123447
$ws.Range("C10").Value = 36.09087621167114
$ws.Range("D10").Value = 0.2597907664266936
$ws.Range("E10").Value = 0.266206582532699
$ws.Range("F10").Value = 0.2600315049119103
$ws.Range("G10").Value = 0.2595427562975338
$ws.Range("H10").Value = 0.260349853614635
$ws.Range("B11").Value = 53.54908979312027
$ws.Range("C11").Value = 45.29978332652574
$ws.Range("D11").Value = 0.2566845647261172
$ws.Range("E11").Value = 0.2652218037878002
$ws.Range("F11").Value = 1.324756094874624
$ws.Range("G11").Value = 0.2563664453588865
$ws.Range("H11").Value = 1.323457669829913
$ws.Range("B12").Value = 71.57082239095257
$ws.Range("C12").Value = 60.17597494014895
$ws.Range("D12").Value = 0.2514405064227141
$ws.Range("E12").Value = 0.2635685322304105
$ws.Range("F12").Value = 4.611752642368373
$ws.Range("G12").Value = 1.344257146668368
$ws.Range("H12").Value = 0.2524342386207764
$ws.Range("B13").Value = 85.69534873469362
$ws.Range("C13").Value = 71.83076988691914
$ws.Range("D13").Value = 2.46775741377846
$ws.Range("E13").Value = 2.355561403309959
$ws.Range("F13").Value = 3.57200642907221
$ws.Range("G13").Value = 0.2466609781383513
$ws.Range("H13").Value = 0.2484028584959173
$ws.Range("B14").Value = 86.94096489842828
$ws.Range("C14").Value = 72.87096263081462
$ws.Range("D14").Value = 1.358565273960585
$ws.Range("E14").Value = 1.309090739428662
$ws.Range("F14").Value = 1.356688655873017
$ws.Range("G14").Value = 1.360463861446458
$ws.Range("H14").Value = 0.2480655756293572
$ws.Range("B15").Value = 85.90817156179996
$ws.Range("C15").Value = 72.05511288773911
$ws.Range("D15").Value = 0.2471398320965782
$ws.Range("E15").Value = 3.403026018938956
$ws.Range("F15").Value = 1.355652091706084
$ws.Range("G15").Value = 1.359354865729184
$ws.Range("H15").Value = 2.458499829307995
$ws.Range("B16").Value = 85.81798607428287
$ws.Range("C16").Value = 71.97702422082635
$ws.Range("D16").Value = 0.247191555749378
$ws.Range("E16").Value = 1.308774817744586
$ws.Range("F16").Value = 1.355473566142384
$ws.Range("G16").Value = 1.359172162998262
$ws.Range("H16").Value = 0.2484051816848972
$ws.Range("B17").Value = 85.7667658756376
$ws.Range("C17").Value = 71.86367437738896
$ws.Range("D17").Value = 0.2471741062084434
$ws.Range("E17").Value = 1.308734845435549
$ws.Range("F17").Value = 3.572262683260283
$ws.Range("G17").Value = 0.2466454528458036
$ws.Range("H17").Value = 2.458273196754498
$ws.Range("B18").Value = 72.23546266634938
$ws.Range("C18").Value = 60.70132432367985
$ws.Range("D18").Value = 0.2512360800664487
$ws.Range("E18").Value = 0.2635087346988858
$ws.Range("F18").Value = 3.523733817105697
$ws.Range("G18").Value = 0.2507973678687586
$ws.Range("H18").Value = 3.516847955431748
$ws.Range("B19").Value = 48.32189163184576
$ws.Range("C19").Value = 40.87286398292552
$ws.Range("D19").Value = 1.321266983515617
$ws.Range("E19").Value = 1.298493216140437
$ws.Range("F19").Value = 1.320402085051094
$ws.Range("G19").Value = 0.2578213462024129
$ws.Range("H19").Value = 4.502444734376246
$ws.Range("B20").Value = 42.14335835588317
$ws.Range("C20").Value = 35.66102027096328
$ws.Range("D20").Value = 2.372678603787931
$ws.Range("E20").Value = 2.32786812018348
$ws.Range("F20").Value = 3.426872045807388
$ws.Range("G20").Value = 1.316852764308057
$ws.Range("H20").Value = 3.423279025068911
$ws.Range("B21").Value = 35.34314762467189
$ws.Range("C21").Value = 30.01787543298084
$ws.Range("D21").Value = 1.310326318843746
$ws.Range("E21").Value = 1.295216677693508
$ws.Range("F21").Value = 1.309731809797891
$ws.Range("G21").Value = 1.310945238605819
$ws.Range("H21").Value = 2.356025539104696
$ws.Range("B22").Value = 32.87240084900301
$ws.Range("C22").Value = 27.93962142466232
$ws.Range("D22").Value = 1.308331030609878
$ws.Range("E22").Value = 1.294606985185827
$ws.Range("F22").Value = 2.353332137422129
$ws.Range("G22").Value = 1.308899503000235
$ws.Range("H22").Value = 1.307068141163787
$ws.Range("B23").Value = 32.06152928340265
$ws.Range("C23").Value = 27.26014945151943
$ws.Range("D23").Value = 1.307703932799422
$ws.Range("E23").Value = 2.322010008390296
$ws.Range("F23").Value = 4.442263711478824
$ws.Range("G23").Value = 0.2623700975436781
$ws.Range("H23").Value = 1.306488059129173
$ws.Range("B24").Value = 29.96576749780159
$ws.Range("C24").Value = 25.49168008102698
$ws.Range("D24").Value = 0.2631221695090668
$ws.Range("E24").Value = 1.293895990711009
$ws.Range("F24").Value = 4.433963649988985
$ws.Range("G24").Value = 1.30651835981251
$ws.Range("H24").Value = 0.263500942548254
$ws.Range("B25").Value = 28.95872638177973
$ws.Range("C25").Value = 24.68398153100106
$ws.Range("D25").Value = 0.2634210020715058
$ws.Range("E25").Value = 0.2673855054977355
$ws.Range("F25").Value = 1.304665320305603
$ws.Range("G25").Value = 1.305615223339886
$ws.Range("H25").Value = 0.2637863201831487
